# Update cryptocurrency price/volume data (refreshed scrape)
# Updated cryptos list on Thu Mar 16 22:49:50 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.970.31'
$ws.Range("E2").Value = '  +1.81%  '
$ws.Range("D3").Value = '1.674.18'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '331.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = '  +1.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.27'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.145'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07151'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9998'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("E13").Value = '  +3.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.49%  '
$ws.Range("D15").Value = '1.666.86'
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.660'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.76%  '
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06548'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9995'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("E20").Value = '  +3.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.917'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.46%  '
$ws.Range("D24").Value = '24.962.21'
$ws.Range("E24").Value = '  +1.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.438'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.391'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("E28").Value = '  +1.53%  '
$ws.Range("D29").Value = '1.853.01'
$ws.Range("E29").Value = '  +0.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.193'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.085'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.798'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08475'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.670'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.56%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("E37").Value = '  -0.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06057'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.230'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02231'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.98%  '
$ws.Range("E41").Value = '  +2.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.246'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9990'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5964'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.27%  '
$ws.Range("E45").Value = '  +8.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.844'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5732'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.58%  '
$ws.Range("E49").Value = '  +1.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07007'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.45%  '
$ws.Range("E51").Value = '  +3.86%  '
